$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Frank Breitenbach" stakeholder row (row 6). Deleting the
# entire row shifts the rows below it (the legend) up by one.
$ws.Rows.Item(6).Delete()

# Match the author's final cursor position after the edit.
$ws.Range("H9").Select() | Out-Null
